# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table (rows 3-8): refreshed Critical Minutes / Good Roaming % numbers ---
$ws.Range("C3").Value = 403
$ws.Range("D3").Value = 83

$ws.Range("C4").Value = 1013
$ws.Range("D4").Value = 92.3

$ws.Range("C5").Value = 6795
$ws.Range("D5").Value = 93

$ws.Range("D6").Value = 97.59999999999999

$ws.Range("C7").Value = 83
$ws.Range("D7").Value = 98

$ws.Range("C8").Value = 8363

# --- "Good Drivers" table (rows 16-24): re-sorted by Driver Vintage (desc) ---
# New newest driver (23.100.0.4, now with an actual sample count + vintage date) is
# inserted at the top; every other row shifts down one slot; two rows pick up a
# Driver Vintage date they were missing before; row 21's vintage date goes from
# blank to 2021-04-27; row 23's vintage date is refreshed to 2020-01-06.

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B16").Value = 445055
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "2024-11-10"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B17").Value = 10661
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "2022-08-29"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B18").Value = 14239
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2022-05-23"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B19").Value = 265400
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").Value = "2022-05-01"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B20").Value = 77849
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "2021-08-18"

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B21").Value = 34244
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "2021-04-27"

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B22").Value = 59673
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = "2020-08-05"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B23").Value = 113652
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "2020-01-06"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B24").Value = 56018
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = "2019-12-14"

Write-Host "Driver summary refreshed for 2025-04-20"
